$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated simulated-game transition probabilities
# (more games simulated => updated probability matrix values)
$ws.Range("B2").Value = 0.2292490118577075
$ws.Range("C2").Value = 0.4703557312252964
$ws.Range("J2").Value = 0.007905138339920948
$ws.Range("P2").Value = 0.1778656126482213
$ws.Range("S2").Value = 0.1146245059288538
$ws.Range("B3").Value = 0.00819672131147541
$ws.Range("C3").Value = 0.01639344262295082
$ws.Range("J3").Value = 0.02459016393442623
$ws.Range("P3").Value = 0.7049180327868853
$ws.Range("S3").Value = 0.2459016393442623
$ws.Range("J4").Value = 0.03703703703703703
$ws.Range("P4").Value = 0.7037037037037037
$ws.Range("S4").Value = 0.2592592592592592
$ws.Range("B6").Value = 0.045662100456621
$ws.Range("D6").Value = 0.0136986301369863
$ws.Range("F6").Value = 0.0410958904109589
$ws.Range("J6").Value = 0.2694063926940639
$ws.Range("O6").Value = 0.0319634703196347
$ws.Range("Q6").Value = 0.1050228310502283
$ws.Range("R6").Value = 0.091324200913242
$ws.Range("S6").Value = 0.4018264840182648
$ws.Range("B7").Value = 0.09444444444444444
$ws.Range("D7").Value = 0.02222222222222222
$ws.Range("F7").Value = 0.05555555555555555
$ws.Range("J7").Value = 0.1833333333333333
$ws.Range("O7").Value = 0.02222222222222222
$ws.Range("Q7").Value = 0.1333333333333333
$ws.Range("R7").Value = 0.1388888888888889
$ws.Range("S7").Value = 0.35
$ws.Range("B8").Value = 0.08488063660477453
$ws.Range("D8").Value = 0.01061007957559682
$ws.Range("F8").Value = 0.07427055702917772
$ws.Range("J8").Value = 0.09549071618037135
$ws.Range("O8").Value = 0.01591511936339523
$ws.Range("Q8").Value = 0.1697612732095491
$ws.Range("R8").Value = 0.143236074270557
$ws.Range("S8").Value = 0.4058355437665783
$ws.Range("B9").Value = 0.09944751381215469
$ws.Range("D9").Value = 0.005524861878453038
$ws.Range("F9").Value = 0.07734806629834254
$ws.Range("J9").Value = 0.1049723756906077
$ws.Range("O9").Value = 0.02762430939226519
$ws.Range("Q9").Value = 0.1491712707182321
$ws.Range("R9").Value = 0.1767955801104972
$ws.Range("S9").Value = 0.3591160220994475
$ws.Range("B10").Value = 0.1130171543895055
$ws.Range("D10").Value = 0.01513622603430878
$ws.Range("F10").Value = 0.09283551967709384
$ws.Range("J10").Value = 0.08577194752774975
$ws.Range("O10").Value = 0.01715438950554995
$ws.Range("Q10").Value = 0.1987891019172553
$ws.Range("R10").Value = 0.08476286579212916
$ws.Range("S10").Value = 0.3925327951564077
$ws.Range("G11").Value = 0.1254612546125461
$ws.Range("J11").Value = 0.09225092250922509
$ws.Range("K11").Value = 0.1845018450184502
$ws.Range("L11").Value = 0.5867158671586716
$ws.Range("S11").Value = 0.01107011070110701
$ws.Range("G12").Value = 0.7267441860465116
$ws.Range("J12").Value = 0.1802325581395349
$ws.Range("K12").Value = 0.005813953488372093
$ws.Range("L12").Value = 0.04651162790697674
$ws.Range("S12").Value = 0.04069767441860465
$ws.Range("G13").Value = 0.7297297297297297
$ws.Range("J13").Value = 0.2162162162162162
$ws.Range("S13").Value = 0.05405405405405406
$ws.Range("F15").Value = 0.01507537688442211
$ws.Range("H15").Value = 0.2110552763819095
$ws.Range("I15").Value = 0.07537688442211055
$ws.Range("J15").Value = 0.3366834170854272
$ws.Range("K15").Value = 0.04522613065326633
$ws.Range("N15").Value = 0.01005025125628141
$ws.Range("O15").Value = 0.05527638190954774
$ws.Range("S15").Value = 0.2512562814070352
$ws.Range("F16").Value = 0.04137931034482759
$ws.Range("H16").Value = 0.1517241379310345
$ws.Range("I16").Value = 0.0896551724137931
$ws.Range("J16").Value = 0.3931034482758621
$ws.Range("K16").Value = 0.1172413793103448
$ws.Range("M16").Value = 0.02758620689655172
$ws.Range("N16").Value = 0.006896551724137931
$ws.Range("O16").Value = 0.06206896551724138
$ws.Range("S16").Value = 0.1103448275862069
$ws.Range("F17").Value = 0.0182370820668693
$ws.Range("H17").Value = 0.2066869300911854
$ws.Range("I17").Value = 0.07598784194528875
$ws.Range("J17").Value = 0.4012158054711246
$ws.Range("K17").Value = 0.1155015197568389
$ws.Range("M17").Value = 0.01519756838905775
$ws.Range("O17").Value = 0.05167173252279635
$ws.Range("S17").Value = 0.1155015197568389
$ws.Range("F18").Value = 0.02314814814814815
$ws.Range("H18").Value = 0.2083333333333333
$ws.Range("I18").Value = 0.08333333333333333
$ws.Range("J18").Value = 0.3518518518518519
$ws.Range("K18").Value = 0.1064814814814815
$ws.Range("M18").Value = 0.02314814814814815
$ws.Range("O18").Value = 0.09259259259259259
$ws.Range("S18").Value = 0.1111111111111111
$ws.Range("F19").Value = 0.0221606648199446
$ws.Range("H19").Value = 0.1939058171745152
$ws.Range("I19").Value = 0.100646352723915
$ws.Range("J19").Value = 0.3481071098799631
$ws.Range("K19").Value = 0.1163434903047091
$ws.Range("M19").Value = 0.02123730378578024
$ws.Range("O19").Value = 0.07571560480147738
$ws.Range("S19").Value = 0.1218836565096953
